$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing data
# columns (B..K, "PercActivations".."totalStd") one to the right (C..L)
# and leaves a fresh, empty column B in place.
$ws.Columns.Item(2).Insert()

# Give the new header cell B1 ("segments") the same bold/bordered header
# style used by the other header cells (copy format from C1, the old B1).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "segments"

# The segment names currently live in column A (rows 2-20). Move them into
# the new column B, and replace column A with a simple 0-based numeric index.
for ($r = 2; $r -le 20; $r++) {
    $segmentName = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $segmentName
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Column B's segment-name cells should use the plain/default style (like the
# other data columns), not the bold bordered style column A uses. Copy the
# default format from one of the untouched data cells (C2) onto B2:B20.
$ws.Range("C2").Copy()
$ws.Range("B2:B20").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
